# Fruta / hortaliza, semanal
# Re-sequence the weekly price rows: the Fecha (D), Volumen (M),
# Precio minimo/maximo/promedio (N/O/P) and Precio $/Kg (S) columns are
# shuffled across rows 2-20 to reflect the corrected weekly ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @(Fecha(serial), Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg)
$rowData = @{
    2  = @(44298, 65, 22000, 22000, 22000, 1100)
    3  = @(44305, 20, 22000, 22000, 22000, 1100)
    4  = @(44307, 30, 22000, 22000, 22000, 1100)
    5  = @(44377, 25, 20000, 20000, 20000, 1000)
    6  = @(44294, 25, 25000, 25000, 25000, 1250)
    7  = @(44403, 50, 20000, 20000, 20000, 1000)
    8  = @(44389, 20, 20000, 20000, 20000, 1000)
    9  = @(44292, 30, 25000, 25000, 25000, 1250)
    10 = @(44406, 20, 20000, 20000, 20000, 1000)
    11 = @(44291, 70, 25000, 25000, 25000, 1250)
    12 = @(44376, 38, 20000, 20000, 20000, 1000)
    13 = @(44300, 45, 22000, 22000, 22000, 1100)
    14 = @(44301, 38, 22000, 22000, 22000, 1100)
    15 = @(44413, 45, 20000, 20000, 20000, 1000)
    16 = @(44382, 24, 20000, 20000, 20000, 1000)
    17 = @(44400, 45, 20000, 20000, 20000, 1000)
    18 = @(44445, 45, 20000, 20000, 20000, 1000)
    19 = @(44385, 36, 20000, 20000, 20000, 1000)
    20 = @(44448, 30, 22000, 22000, 22000, 1100)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[5]   # S - Precio $/Kg
}
